# "plot week by week" — re-point the task's ParameterNameList workbook at
# the new per-week OIS parameter set and drop the now-unused "all (2)" sheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the "all (2)" sheet entirely (it duplicated a subset of "all").
$oldSheet = $wb.Worksheets.Item("all (2)")
[void]$oldSheet.Delete()

$ws = $wb.Worksheets.Item("all")

# Row 1 header stays the same (plist / nlist / typelist).
$ws.Range("A1").Value = "plist"
$ws.Range("B1").Value = "nlist"
$ws.Range("C1").Value = "typelist"

# Row 2: immediate object accuracy — only the plist name changes case.
$ws.Range("A2").Value = "ois_ImmediateObjectAccuracy"
$ws.Range("B2").Value = "Immediate object accuracy (%)"
$ws.Range("C2").Value = "num"

# Row 3: now semantic accuracy (immediate).
$ws.Range("A3").Value = "ois_ImmediateSemanticAccuracy"
$ws.Range("B3").Value = "Delayed object accuracy (%)"
$ws.Range("C3").Value = "num"

# Row 4: now immediate location error.
$ws.Range("A4").Value = "ois_ImmediateLocationError"
$ws.Range("B4").Value = "Immediate location error (cm)"
$ws.Range("C4").Value = "num"

# Row 5: now delayed object accuracy.
$ws.Range("A5").Value = "ois_DelayedObjectAccuracy"
$ws.Range("B5").Value = "Delayed location error (cm)"
$ws.Range("C5").Value = "num"

# Row 6: now delayed semantic accuracy.
$ws.Range("A6").Value = "ois_DelayedSemanticAccuracy"
$ws.Range("B6").Value = "Immediate semantic accuracy (%)"
$ws.Range("C6").Value = "num"

# Row 7: now delayed location error.
$ws.Range("A7").Value = "ois_DelayedLocationError"
$ws.Range("B7").Value = "Delayed semantic accuracy (%)"
$ws.Range("C7").Value = "num"

# Rows 8-11: new reaction-time parameters, no nlist label.
$ws.Range("A8").Value = "ois_rt_identification_immediate"
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = "num"

$ws.Range("A9").Value = "ois_rt_localisation_immediate"
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value = "num"

$ws.Range("A10").Value = "ois_rt_identification_delayed"
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = "num"

$ws.Range("A11").Value = "ois_rt_localisation_delayed"
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = "num"

# Rows 12-13 no longer exist.
$ws.Range("A12:C13").ClearContents()

[void]$ws.Activate()
[void]$ws.Range("A2").Select()
